# Natmi following Dr Hou advice
# Adds "ECs" as a third sending/target cluster alongside the existing
# "FAPs" and "sCs" clusters for the Comp -> Itga5 ligand/receptor pair,
# expanding the 2x3 combination grid (rows 2-7) into a full 3x3 grid
# (rows 2-10) and refreshing the associated NATMI statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @(
    @("ECs", "Comp", "Itga5", "ECs", 2, 0.6666666666666666, 0.569336, 1.708008, 0.01936485544401529, 0.01936485544401529, 3, 1, 34.07074633333333, 102.212239, 0.5171464495142372, 0.5171464495142373, 19.39770243443466, 174.579321909912, 0.01001446623822896, 0.01001446623822896),
    @("ECs", "Comp", "Itga5", "FAPs", 2, 0.6666666666666666, 0.569336, 1.708008, 0.01936485544401529, 0.01936485544401529, 3, 1, 27.685497, 83.056491, 0.420227262899125, 0.4202272628991251, 15.762350119992, 141.861151079928, 0.008137640199675767, 0.008137640199675769),
    @("ECs", "Comp", "Itga5", "sCs", 2, 0.6666666666666666, 0.569336, 1.708008, 0.01936485544401529, 0.01936485544401529, 3, 1, 4.125957666666666, 12.377873, 0.06262628758663766, 0.06262628758663766, 2.349056234109333, 21.141506106984, 0.001212749006110568, 0.001212749006110568),
    @("FAPs", "Comp", "Itga5", "ECs", 3, 1, 28.31465866666666, 84.94397599999999, 0.9630679809930072, 0.9630679809930072, 3, 1, 34.07074633333333, 102.212239, 0.5171464495142372, 0.5171464495142373, 964.7015529469179, 8682.313976522262, 0.4980471870113785, 0.4980471870113786),
    @("FAPs", "Comp", "Itga5", "FAPs", 3, 1, 28.31465866666666, 84.94397599999999, 0.9630679809930072, 0.9630679809930072, 3, 1, 27.685497, 83.056491, 0.420227262899125, 0.4202272628991251, 783.9053975720238, 7055.148578148215, 0.404707421638478, 0.4047074216384781),
    @("FAPs", "Comp", "Itga5", "sCs", 3, 1, 28.31465866666666, 84.94397599999999, 0.9630679809930072, 0.9630679809930072, 3, 1, 4.125957666666666, 12.377873, 0.06262628758663766, 0.06262628758663766, 116.8250830047831, 1051.425747043048, 0.06031337234315055, 0.06031337234315055),
    @("sCs", "Comp", "Itga5", "ECs", 3, 1, 0.516483, 1.549449, 0.01756716356297749, 0.01756716356297749, 3, 1, 34.07074633333333, 102.212239, 0.5171464495142372, 0.5171464495142373, 17.596961278479, 158.372651506311, 0.009084796264629684, 0.009084796264629685),
    @("sCs", "Comp", "Itga5", "FAPs", 3, 1, 0.516483, 1.549449, 0.01756716356297749, 0.01756716356297749, 3, 1, 27.685497, 83.056491, 0.420227262899125, 0.4202272628991251, 14.299088547051, 128.691796923459, 0.007382201060971271, 0.007382201060971272),
    @("sCs", "Comp", "Itga5", "sCs", 3, 1, 0.516483, 1.549449, 0.01756716356297749, 0.01756716356297749, 3, 1, 4.125957666666666, 12.377873, 0.06262628758663766, 0.06262628758663766, 2.130986993553, 19.178882941977, 0.00110016623737653, 0.00110016623737653)
)

$rowIndex = 2
foreach ($row in $data) {
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "$($columns[$i])$rowIndex"
        $ws.Range($cellRef).Value = $row[$i]
    }
    $rowIndex++
}
